# Applies the NIT-9011461914 "Estado de Cuenta" update:
#  - Reorders the debtor rows: BLAS JOSE HERNANDEZ GENES' periods now occupy
#    rows 16-22 (most recent period 2401 first, oldest 2307 last), and
#    LIDIS DEL CARMEN ALVAREZ CAMARGO's single period (2201) moves to row 23.
#  - Column widths widen slightly to fit the refreshed "best fit" content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data: Doc No. (C), Name (D), Period (E), Valor Mora (F), Salario Basico (G) ---
$rows = @(
    @{ Row = 16; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2401"; F = 35112; G = 1160000 },
    @{ Row = 17; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2312"; F = 35112; G = 1160000 },
    @{ Row = 18; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2311"; F = 35112; G = 1160000 },
    @{ Row = 19; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2310"; F = 35112; G = 1160000 },
    @{ Row = 20; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2309"; F = 35112; G = 1160000 },
    @{ Row = 21; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2308"; F = 46400; G = 1160000 },
    @{ Row = 22; Doc = "15025168"; Name = "BLAS JOSE HERNANDEZ GENES";           Period = "2307"; F = 46400; G = 1160000 },
    @{ Row = 23; Doc = "45487733"; Name = "LIDIS DEL CARMEN ALVAREZ CAMARGO";    Period = "2201"; F = 40000; G = 1000000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Doc
    $ws.Range("D$n").Value = $r.Name
    $ws.Range("E$n").Value = $r.Period
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}

# --- Column widths (refreshed "best fit" sizing) ---
$ws.Columns("B").ColumnWidth = 17.709635416666668
$ws.Columns("C").ColumnWidth = 15.893229166666666
$ws.Columns("E").ColumnWidth = 12.709635416666666
$ws.Columns("F").ColumnWidth = 9.346354166666666
$ws.Columns("G").ColumnWidth = 13.529947916666666
$ws.Columns("H").ColumnWidth = 18.529947916666668
$ws.Columns("I").ColumnWidth = 17.256510416666668
$ws.Columns("J").ColumnWidth = 14.166666666666666
